$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 6 (Prof. Edward Casper IV / rene57@example.net), shifting rows up
$ws.Rows.Item(6).Delete()

# Set active cell selection to A6 (where cursor lands after the delete)
$ws.Range("A6").Select() | Out-Null
